$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_parameter")

# Switch solver from Gurobi.jl to HiGHS.jl for both the db_lp_solver and
# db_mip_solver parameters (row 2, columns E and F).
$ws.Range("E2").Value = "HiGHS.jl"
$ws.Range("F2").Value = "HiGHS.jl"

# Reflect the cursor/selection left on the active sheet when the workbook
# was last saved.
$ws.Activate() | Out-Null
$ws.Range("F8").Select() | Out-Null
